$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (also reflected in workbook.xml <sheet name=.../>)
$ws.Name = "Through 2022-08-09"

# Update the header label for the August column
$ws.Range("A9").Value = "August (through 08-09)"

# Update August row (row 9) values for columns B..I
$ws.Range("B9").Value = 10
$ws.Range("C9").Value = 15
$ws.Range("D9").Value = 23
$ws.Range("E9").Value = 17
$ws.Range("F9").Value = 13
$ws.Range("G9").Value = 52
$ws.Range("H9").Value = 60
$ws.Range("I9").Value = 50

# Update Total row (row 10) values for columns B..I
$ws.Range("B10").Value = 172
$ws.Range("C10").Value = 317
$ws.Range("D10").Value = 488
$ws.Range("E10").Value = 442
$ws.Range("F10").Value = 317
$ws.Range("G10").Value = 673
$ws.Range("H10").Value = 970
$ws.Range("I10").Value = 1020
